# BOM update: add R8 "Potentiometer" line, shrink the R2/R3/R5/R8/R10 group
# down to R2/R3/R5 (R8 now its own pot, R10 removed), widen the Description
# column, reset row 10's height, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: "R2, R3, R5, R8, R10" group loses R8 & R10 -> "R2, R3, R5" -----
$ws.Range("A5").Value = "R2, R3, R5"
$ws.Range("F5").Value = 3
$ws.Range("I5").Value = 3 * 0.0037

# --- New row 13: R8 is now its own Potentiometer line ---------------------
$ws.Cells.Item(13, 1).Value = "R8"
$ws.Cells.Item(13, 2).Value = "Potentiometer"
$ws.Cells.Item(13, 3).Value = "TRIMMER 10K OHM 0.25W SMD"
$ws.Cells.Item(13, 4).Value = "3314Z-1-103E"
$ws.Cells.Item(13, 5).Value = "3314Z-103ECT-ND"
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 2.006
$ws.Cells.Item(13, 8).Value = "10"
$ws.Cells.Item(13, 9).Value = 2.006
$ws.Cells.Item(13, 10).Value = "http://www.digikey.com/product-detail/en/bourns-inc/3314Z-1-103E/3314Z-103ECT-ND/253541"

# C13/D13 use a small Arial font (pasted-in datasheet style) instead of the
# sheet's normal Calibri; D13 keeps the wrap/vcenter treatment C10 also uses.
$rowFont = $ws.Range("C13:D13").Font
$rowFont.Name = "Arial"
$rowFont.Size = 9
$rowFont.Color = 0
$ws.Range("D13").WrapText = $true
$ws.Range("D13").VerticalAlignment = -4108

# Row 10 no longer needs the taller auto-wrap height it had before.
$ws.Rows(10).RowHeight = 15.75

# --- Column C is a bit wider now, and no longer a "best fit" column -------
$ws.Columns("C").ColumnWidth = 40

# --- Active cell moves from the old edit point to I6 -----------------------
$ws.Range("I6").Select()
